$d = $word.ActiveDocument

# 1. Merge the split "Nuke"/"Maya" runs into a single run of plain text.
$d.Content.Find.Execute(
    "I would like to give two examples for this. One for Nuke and other for Maya",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I would like to give two examples for this. One for Nuke and other for Maya",
    2)

# 2. Remove the leading "EC2" paragraph and the empty paragraph after it
#    (the first two paragraphs of the document).
$r = $d.Range($d.Paragraphs(1).Range.Start, $d.Paragraphs(2).Range.End)
$r.Delete()
